$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.742.25'
$ws.Range("D3").Value = '1.758.28'
$ws.Range("E3").Value = '  -2.07%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '325.78'
$ws.Range("E5").Value = '  -3.32%  '
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("D7").Value = '0.4409'
$ws.Range("E7").Value = '  -2.87%  '
$ws.Range("D8").Value = '0.3727'
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").Value = '45.42'
$ws.Range("E9").Value = '  +0.63%  '
$ws.Range("D10").Value = '0.07533'
$ws.Range("E10").Value = '  -0.80%  '
$ws.Range("D11").Value = '1.125'
$ws.Range("E11").Value = '  -1.34%  '
$ws.Range("E12").Value = '  -0.27%  '
$ws.Range("E13").Value = '  -2.13%  '
$ws.Range("D14").Value = '6.207'
$ws.Range("E14").Value = '  -0.99%  '
$ws.Range("D15").Value = '7.409'
$ws.Range("E15").Value = '  -0.64%  '
$ws.Range("D16").Value = '1.758.96'
$ws.Range("E16").Value = '  -2.13%  '
$ws.Range("D17").Value = '0.00001072'
$ws.Range("E17").Value = '  -1.28%  '
$ws.Range("D18").Value = '88.29'
$ws.Range("E18").Value = '  +8.87%  '
$ws.Range("D19").Value = '0.06211'
$ws.Range("E19").Value = '  -7.84%  '
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("D21").Value = '17.39'
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").Value = '6.182'
$ws.Range("E22").Value = '  -2.90%  '
$ws.Range("D23").Value = '0.5320'
$ws.Range("E23").Value = '  -3.51%  '
$ws.Range("D24").Value = '27.769.99'
$ws.Range("E24").Value = '  -1.59%  '
$ws.Range("D25").Value = '11.65'
$ws.Range("E25").Value = '  -0.93%  '
$ws.Range("D26").Value = '2.315'
$ws.Range("E26").Value = '  -4.50%  '
$ws.Range("D27").Value = '20.65'
$ws.Range("E27").Value = '  +0.72%  '
$ws.Range("D28").Value = '153.12'
$ws.Range("E28").Value = '  +0.31%  '
$ws.Range("D29").Value = '2.371'
$ws.Range("E29").Value = '  +1.21%  '
$ws.Range("D30").Value = '1.957.65'
$ws.Range("E30").Value = '  -2.28%  '
$ws.Range("D31").Value = '128.30'
$ws.Range("E31").Value = '  -3.44%  '
$ws.Range("D32").Value = '1.219'
$ws.Range("E32").Value = '  -0.60%  '
$ws.Range("D33").Value = '0.09363'
$ws.Range("E33").Value = '  -0.57%  '
$ws.Range("D34").Value = '5.755'
$ws.Range("E34").Value = '  -0.42%  '
$ws.Range("D35").Value = '3.650'
$ws.Range("E35").Value = '  -9.62%  '
$ws.Range("D36").Value = '12.71'
$ws.Range("E36").Value = '  +5.85%  '
$ws.Range("D37").Value = '0.02334'
$ws.Range("E37").Value = '  -0.14%  '
$ws.Range("E38").Value = '  -6.98%  '
$ws.Range("D39").Value = '0.06138'
$ws.Range("E39").Value = '  -2.66%  '
$ws.Range("D40").Value = '0.6495'
$ws.Range("E40").Value = '  -0.96%  '
$ws.Range("D41").Value = '5.084'
$ws.Range("E41").Value = '  -2.22%  '
$ws.Range("E42").Value = '  -0.43%  '
$ws.Range("D43").Value = '7.998'
$ws.Range("E43").Value = '  -3.89%  '
$ws.Range("D44").Value = '1.419'
$ws.Range("E44").Value = '  -4.21%  '
$ws.Range("D45").Value = '1.000'
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("D46").Value = '13.80'
$ws.Range("E46").Value = '  -2.27%  '
$ws.Range("D47").Value = '0.6003'
$ws.Range("E47").Value = '  -1.20%  '
$ws.Range("D48").Value = '3.754'
$ws.Range("E48").Value = '  -2.35%  '
$ws.Range("D49").Value = '126.50'
$ws.Range("E49").Value = '  -2.52%  '
$ws.Range("D50").Value = '1.987'
$ws.Range("E50").Value = '  -1.52%  '
$ws.Range("E51").Value = '  -3.08%  '
